{"js": "// Sequential replacement of the date line and every division-problem\n// cell in the table, in document order. Several source values repeat\n// (e.g. \"44\u00f72=\" and \"18\u00f75=\" each occur twice) but map to different\n// targets depending on position, so we must match by order rather than\n// by unique text.\nconst replacements = [\n  \"2025-08-26 Tuesday\",\n  \"18\u00f75=\",\n  \"17\u00f73=\",\n  \"73\u00f78=\",\n  \"74\u00f75=\",\n  \"12\u00f72=\",\n  \"14\u00f74=\",\n  \"57\u00f79=\",\n  \"83\u00f79=\",\n  \"12\u00f74=\",\n  \"94\u00f72=\",\n  \"32\u00f74=\",\n  \"85\u00f73=\",\n  \"29\u00f73=\",\n  \"43\u00f72=\",\n  \"80\u00f77=\",\n  \"27\u00f79=\",\n  \"22\u00f74=\",\n  \"46\u00f79=\",\n  \"16\u00f76=\",\n  \"38\u00f72=\",\n  \"10\u00f79=\",\n  \"83\u00f76=\",\n  \"42\u00f76=\",\n  \"93\u00f73=\",\n  \"79\u00f79=\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"\") continue;\n  if (idx >= replacements.length) break;\n  para.insertText(replacements[idx], \"Replace\");\n  idx++;\n}\n\nawait context.sync();\n", "ps1": "# Sequential replacement of the date line and every division-problem\n# cell in the table, in document order. Several source values repeat\n# (e.g. \"44\u00f72=\" and \"18\u00f75=\" each occur twice) but map to different\n# targets depending on position, so we match by order rather than by\n# unique text.\n$replacements = @(\n  \"2025-08-26 Tuesday\",\n  \"18\u00f75=\",\n  \"17\u00f73=\",\n  \"73\u00f78=\",\n  \"74\u00f75=\",\n  \"12\u00f72=\",\n  \"14\u00f74=\",\n  \"57\u00f79=\",\n  \"83\u00f79=\",\n  \"12\u00f74=\",\n  \"94\u00f72=\",\n  \"32\u00f74=\",\n  \"85\u00f73=\",\n  \"29\u00f73=\",\n  \"43\u00f72=\",\n  \"80\u00f77=\",\n  \"27\u00f79=\",\n  \"22\u00f74=\",\n  \"46\u00f79=\",\n  \"16\u00f76=\",\n  \"38\u00f72=\",\n  \"10\u00f79=\",\n  \"83\u00f76=\",\n  \"42\u00f76=\",\n  \"93\u00f73=\",\n  \"79\u00f79=\"\n)\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n$idx = 0\n\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs($i)\n  $t = $p.Range.Text\n  $t2 = $t.TrimEnd([char]13, [char]7)\n  if ($t2.Length -gt 0) {\n    if ($idx -lt $replacements.Length) {\n      $p.Range.Text = $replacements[$idx]\n    }\n    $idx = $idx + 1\n  }\n}\n"}
